$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# Insert the two new sheets in order right after LoginData:
# LoginData, SearchProduct, CheckOut
$wsSearch = $wb.Worksheets.Add($null, $loginSheet)
$wsSearch.Name = "SearchProduct"
$wsCheckout = $wb.Worksheets.Add($null, $wsSearch)
$wsCheckout.Name = "CheckOut"

# --- SearchProduct sheet data ---
$wsSearch.Range("A1:A2").NumberFormat = "@"
$wsSearch.Range("A1").Value = "Product"
$wsSearch.Range("A2").Value = "Sauce Labs Bike Light"
[void]$wsSearch.Range("A2").Select()

# --- CheckOut sheet data ---
$wsCheckout.Range("A1:C3").NumberFormat = "@"
$wsCheckout.Range("A1").Value = "FirstName"
$wsCheckout.Range("B1").Value = "LastName"
$wsCheckout.Range("C1").Value = "Zipcode"
$wsCheckout.Range("A2").Value = "Henry"
$wsCheckout.Range("B2").Value = "Santa"
$wsCheckout.Range("C2").Value = "700156"
[void]$wsCheckout.Range("K4").Select()

# CheckOut is the tab that ends up selected/active when the workbook is saved
[void]$wsCheckout.Activate()
